$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update clinic name labels (shorten full names to abbreviations)
$ws.Range("A2").Value = "CHCF"
$ws.Range("A3").Value = "THS"
$ws.Range("A4").Value = "PPH"

# Update D3:E4 values from 1 to 3
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 3

# Update row heights: text is shorter now (no wrap needed), rows shrink to 17pt
$ws.Rows.Item(2).RowHeight = 17
$ws.Rows.Item(3).RowHeight = 17
$ws.Rows.Item(4).RowHeight = 17

# Update selection to A5
$ws.Range("A5").Select()
